$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Move the exercise link from row 4 (session 3) to row 5 (session 4),
# and update the filename from e03 to e04.
$ws.Range("F4").ClearContents()
$ws.Range("F5").Value = "exercises/e04.html"

# Update the selection to reflect where the user last clicked (F6).
$ws.Range("F6").Select()
